$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the coin table with refreshed prices / 1h volume percentages.
# Column D holds price text; where the new value is a plain decimal number
# (e.g. "212.74"), force the cell to Text format first so Excel keeps it as
# a string instead of silently converting it to a numeric value (matching
# the original inline-string storage of these cells).

# Row 2
$ws.Range("D2").Value = "26.060.95"
$ws.Range("E2").Value = "  +3.09%  "

# Row 3
$ws.Range("D3").Value = "1.602.63"
$ws.Range("E3").Value = "  +3.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.74"
$ws.Range("E5").Value = "  +2.67%  "

# Row 6
$ws.Range("E6").Value = "  -0.08%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.486"
$ws.Range("E7").Value = "  +1.59%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.247"
$ws.Range("E8").Value = "  +1.93%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0613"
$ws.Range("E9").Value = "  +0.53%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.04"
$ws.Range("E10").Value = "  +1.67%  "

# Row 11
$ws.Range("E11").Value = "  +4.72%  "

# Row 12
$ws.Range("D12").Value = "1.827.15"
$ws.Range("E12").Value = "  +3.36%  "

# Row 13
$ws.Range("D13").Value = "1.603.69"
$ws.Range("E13").Value = "  +3.57%  "

# Row 14
$ws.Range("E14").Value = "  +0.28%  "

# Row 15
$ws.Range("E15").Value = "  +1.21%  "

# Row 16
$ws.Range("D16").Value = "26.070.03"
$ws.Range("E16").Value = "  +3.21%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.35"
$ws.Range("E17").Value = "  +2.75%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0720"
$ws.Range("E18").Value = "  +1.56%  "

# Row 19
$ws.Range("E19").Value = "  -0.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.89"
$ws.Range("E20").Value = "  +8.81%  "

# Row 21
$ws.Range("E21").Value = "  +2.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.28"
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.98"
$ws.Range("E23").Value = "  +2.34%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.85"
$ws.Range("E24").Value = "  +13.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.44"
$ws.Range("E25").Value = "  +1.14%  "

# Row 26
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("E27").Value = "  -4.83%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.16"
$ws.Range("E28").Value = "  +2.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.42"
$ws.Range("E29").Value = "  +0.34%  "

# Row 30
$ws.Range("E30").Value = "  +1.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("E31").Value = "  +0.97%  "

# Row 32
$ws.Range("E32").Value = "  +2.64%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.95"
$ws.Range("E33").Value = "  -0.47%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("E34").Value = "  +1.59%  "

# Row 35
$ws.Range("E35").Value = "  +0.92%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0165"
$ws.Range("E36").Value = "  +10.75%  "

# Row 37
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.123.77"
$ws.Range("E37").Value = "  +3.50%  "

# Row 38
$ws.Range("E38").Value = "  -0.03%  "

# Row 39
$ws.Range("E39").Value = "  +2.59%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.781"
$ws.Range("E40").Value = "  +1.89%  "

# Row 41
$ws.Range("E41").Value = "  -0.60%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.782"
$ws.Range("E42").Value = "  -2.40%  "

# Row 43
$ws.Range("E43").Value = "  +2.29%  "

# Row 44
$ws.Range("D44").Value = "1.740.13"
$ws.Range("E44").Value = "  +3.42%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.74"
$ws.Range("E45").Value = "  +0.15%  "

# Row 46
$ws.Range("E46").Value = "  +3.94%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.43"
$ws.Range("E47").Value = "  +2.08%  "

# Row 48
$ws.Range("E48").Value = "  +0.41%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.409"
$ws.Range("E49").Value = "  +1.24%  "

# Row 50
$ws.Range("E50").Value = "  +0.08%  "

# Row 51
$ws.Range("D51").Value = "0.0₇0925"
$ws.Range("E51").Value = "  -16.84%  "
